$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B16 was an inline string "3"; change it to a real number 3
$ws.Range("B16").Value = 3

# Add new row 17 with the new annotation data
$ws.Range("A17").Value = "Ruilin"
$ws.Range("B17").Value = "'4"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "can solve more complex tasks"
$ws.Range("D17").Value = "APC"
$ws.Range("E17").Value = "RES"
$ws.Range("F17").Value = "9cb2103f-10a8-4188-b35f-b6e342d90889"
$ws.Range("G17").Value = "rJwelMbR-_annotated.xlsx"
$ws.Range("H17").Value = "The authors show through several experiments that the divide and conquer (DnC) technique can solve more complex tasks than can be solved with conventional policy gradient methods (TRPO is used as the baseline)."
